$wb = $excel.ActiveWorkbook
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Sheet2"
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
Write-Host "Last sheet before move: $($last.Name)"
$newSheet.Move($null, $wb.Worksheets.Item("Sheet1"))
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
